# Apply the crypto price/volume updates for this commit.
# The Price column (D) holds numeric-looking text (e.g. "1.003") that must
# stay text, not get auto-converted to a number by Excel, so force the
# cell format to Text ("@") before writing the value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.570.99"
$ws.Range("E2").Value = "  -0.81%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.883.04"
$ws.Range("E3").Value = "  -0.76%  "

$ws.Range("E4").Value = "  +0.52%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.94"
$ws.Range("E5").Value = "  -4.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.46%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4857"
$ws.Range("E7").Value = "  -1.58%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2887"
$ws.Range("E8").Value = "  -2.33%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06607"
$ws.Range("E9").Value = "  -1.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.889.53"
$ws.Range("E10").Value = "  -0.24%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.71"
$ws.Range("E11").Value = "  -0.30%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07208"
$ws.Range("E12").Value = "  -0.50%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "88.34"
$ws.Range("E13").Value = "  +0.42%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.993"
$ws.Range("E14").Value = "  -1.12%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6637"
$ws.Range("E15").Value = "  -1.93%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.546.16"
$ws.Range("E16").Value = "  -0.53%  "

$ws.Range("B17").Value = "Dai"
$ws.Range("C17").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.004"
$ws.Range("E17").Value = "  +0.64%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007824"
$ws.Range("E18").Value = "  -1.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.97"
$ws.Range("E19").Value = "  -0.24%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.133.72"
$ws.Range("E20").Value = "  +0.18%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.004"
$ws.Range("E21").Value = "  +0.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.747"
$ws.Range("E22").Value = "  -1.51%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "188.19"
$ws.Range("E23").Value = "  +16.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.024"
$ws.Range("E24").Value = "  +1.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.231"
$ws.Range("E25").Value = "  -0.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.05"
$ws.Range("E26").Value = "  +2.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.50"
$ws.Range("E27").Value = "  +6.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.845"
$ws.Range("E28").Value = "  -4.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.416"
$ws.Range("E29").Value = "  -0.78%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.240"
$ws.Range("E30").Value = "  -0.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08986"
$ws.Range("E31").Value = "  +1.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.909"
$ws.Range("E32").Value = "  -2.57%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05238"
$ws.Range("E33").Value = "  -0.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7299"
$ws.Range("E34").Value = "  -1.24%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.075"
$ws.Range("E35").Value = "  -4.70%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.717"
$ws.Range("E36").Value = "  +2.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01819"
$ws.Range("E37").Value = "  -1.66%  "

$ws.Range("E38").Value = "  -1.59%  "

$ws.Range("E39").Value = "  -2.87%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.057"
$ws.Range("E40").Value = "  -6.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4322"
$ws.Range("E41").Value = "  +0.71%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "103.92"
$ws.Range("E42").Value = "  -0.56%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9988"
$ws.Range("E43").Value = "  -0.15%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.590"
$ws.Range("E44").Value = "  -4.53%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1328"
$ws.Range("E45").Value = "  +1.78%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.303"
$ws.Range("E46").Value = "  -3.56%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05826"
$ws.Range("E47").Value = "  +0.47%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.654"
$ws.Range("E48").Value = "  +3.02%  "

$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.31"
$ws.Range("E49").Value = "  +0.89%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.402"
$ws.Range("E50").Value = "  +3.15%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3885"

